# Updated symbol list on Wed Jan 11 16:53:25 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) values for the
# crypto listing on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "277.04";     E = "0.69%" }
    @{ Row = 3;  D = "27.26";      E = "0.41%" }
    @{ Row = 4;  D = "4.836";      E = "1.48%" }
    @{ Row = 5;  D = "0.06329";    E = "0.25%" }
    @{ Row = 6;  D = "7.018";      E = "1.20%" }
    @{ Row = 7;  D = "1.377";      E = "2.96%" }
    @{ Row = 8;  D = "0.8903";     E = "1.46%" }
    @{ Row = 9;  D = "0.1518";     E = "0.13%" }
    @{ Row = 10; D = "0.05276";    E = "4.94%" }
    @{ Row = 11; D = "0.07447";    E = "0.38%" }
    @{ Row = 12; D = "0.02895";    E = "0.98%" }
    @{ Row = 13; D = "0.08955";    E = "-0.77%" }
    @{ Row = 14; D = "0.001564";   E = "-0.97%" }
    @{ Row = 15; D = "0.0006358";  E = "0.13%" }
    @{ Row = 16; D = "0.006040";   E = "3.45%" }
    @{ Row = 17; D = "3.471";      E = "0.66%" }
    @{ Row = 18; D = "3.298";      E = "-0.16%" }
    @{ Row = 19; D = "2.234";      E = "-1.68%" }
    @{ Row = 21; E = "0.57%" }
    @{ Row = 22; D = "3.909";      E = "-0.05%" }
    @{ Row = 23; D = "0.1507";     E = "9.18%" }
    @{ Row = 24; D = "0.04380";    E = "-0.71%" }
    @{ Row = 25; D = "0.001178";   E = "0.26%" }
    @{ Row = 26; D = "0.004242";   E = "10.57%" }
    @{ Row = 28; E = "-1.75%" }
    @{ Row = 29; E = "-14.98%" }
    @{ Row = 40; D = "0.03974";    E = "-2.84%" }
    @{ Row = 41; D = "0.006661";   E = "-2.41%" }
    @{ Row = 42; D = "0.1413";     E = "20.55%" }
    @{ Row = 43; D = "0.001981";   E = "-10.89%" }
    @{ Row = 44; D = "0.01173";    E = "1.65%" }
    @{ Row = 45; D = "0.00005329"; E = "2.75%" }
    @{ Row = 46; D = "1.561";      E = "4.79%" }
    @{ Row = 47; D = "0.01850";    E = "-19.69%" }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        # Leading apostrophe forces Excel to store the numeric-looking text
        # ("277.04", "0.06329", ...) as a literal text value, matching the
        # existing inline-string cell type instead of auto-converting to a
        # number.
        $ws.Cells.Item($u.Row, 4).Value = "'" + $u.D
    }
    if ($u.ContainsKey("E")) {
        # Same trick for the percentage-formatted text in column E.
        $ws.Cells.Item($u.Row, 5).Value = "'" + $u.E
    }
}
